# Weekly data refresh: a new week's price observation for "Frambuesa"
# (raspberry) at Mercado Mayorista Lo Valledor de Santiago is prepended
# to the existing data block, pushing the previously-existing rows
# 172:205 down to 173:206 (dimension grows from A1:T205 to A1:T206).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row above the current row 172; Excel shifts every
# row at/after 172 down by one (so old 172 -> 173, ..., old 205 -> 206)
# and carries the row-172 number formatting (date style on column D) down
# with it, matching how the existing rows are formatted.
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new week's observation.
$ws.Range("A172").Value = 6
$ws.Range("B172").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C172").Value = "Metropolitana"
$ws.Range("D172").Value = 44694
$ws.Range("E172").Value = 13
$ws.Range("F172").Value = "Fruta"
$ws.Range("G172").Value = 100101
$ws.Range("H172").Value = "Berries"
$ws.Range("I172").Value = 100101004
$ws.Range("J172").Value = "Frambuesa"
$ws.Range("K172").Value = "Sin especificar"
$ws.Range("L172").Value = "Primera"
$ws.Range("M172").Value = 75
$ws.Range("N172").Value = 12000
$ws.Range("O172").Value = 12000
$ws.Range("P172").Value = 12000
$ws.Range("Q172").Value = "`$/bandeja 2 kilos"
$ws.Range("R172").Value = "Provincia de Curicó"
$ws.Range("S172").Value = 6000
$ws.Range("T172").Value = 2
